$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I and J ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-17: new numeric columns I and J ---
$data = @(
  @(1, 3),
  @(2, 4),
  @(5, 9),
  @(8, 8),
  @(1, 5),
  @(1, 6),
  @(8, 9),
  @(1, 5),
  @(1, 4),
  @(1, 6),
  @(1, 4),
  @(1, 4),
  @(6, 7),
  @(1, 3),
  @(5, 7),
  @(5, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 9).Value = $data[$i][0]
  $ws.Cells.Item($r, 10).Value = $data[$i][1]
}
